# Add newly collected submission rows (79-90) to the first worksheet
# ("八位序列号收集收集结果yd5"), matching the appended rows in the
# upstream "收集结果" workbook update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dateFormat = "yyyy/m/d h:mm:ss;@"

$rows = @(
    @{ Row = 79;  A = "陆十柒（羽翼春秋）";        B = 45954.9465046296; C = "3d66f1c1"; D = "2425167078" },
    @{ Row = 80;  A = "大声发";                    B = 45955.4026851852; C = "6e70d1ea"; D = "176324771" },
    @{ Row = 81;  A = "盛师傅";                    B = 45955.4330092593; C = "8a0ddaac"; D = "1470740944" },
    @{ Row = 82;  A = "江東.";                     B = 45955.8014467593; C = "35557dba"; D = "1508574214" },
    @{ Row = 83;  A = "℘̶敗̶給̶鐘̶意";                B = 45955.8967476852; C = "e3a9c24f"; D = "1718358607" },
    @{ Row = 84;  A = "🐏如画";                    B = 45957.8094444445; C = "ec74cad7"; D = "1828429965" },
    @{ Row = 85;  A = "ly";                        B = 45957.9504050926; C = "7d895744"; D = "1542303450" },
    @{ Row = 86;  A = "L";                         B = 45959.1914583333; C = "60527ee0"; D = "3579627143" },
    @{ Row = 87;  A = "终若.";                     B = 45959.3531481482; C = "f738671."; D = "391111874" },
    @{ Row = 88;  A = "哭";                        B = 45960.4028356482; C = "60527ee0"; D = "2564982476" },
    @{ Row = 89;  A = "Voyager 1";                 B = 45961.5695833333; C = "c720d4f6"; D = "1508392412" },
    @{ Row = 90;  A = "shmily";                    B = 45961.5820138889; C = "f25dfbe5"; D = "3202344808" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 2).NumberFormat = $dateFormat
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
